$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 163, shifting existing rows 163-173 down to 164-174
$ws.Rows.Item(163).Insert()

# Populate the new row 163 with the new weekly price-report entry
$ws.Cells.Item(163, 1).Value  = 11
$ws.Cells.Item(163, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(163, 3).Value  = "Bíobío"
$ws.Cells.Item(163, 4).Value  = 44714
$ws.Cells.Item(163, 5).Value  = 8
$ws.Cells.Item(163, 6).Value  = "Fruta"
$ws.Cells.Item(163, 7).Value  = 100108
$ws.Cells.Item(163, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(163, 9).Value  = 100108005
$ws.Cells.Item(163, 10).Value = "Piña"
$ws.Cells.Item(163, 11).Value = "Caramelo"
$ws.Cells.Item(163, 12).Value = "Segunda"
$ws.Cells.Item(163, 13).Value = 200
$ws.Cells.Item(163, 14).Value = 17000
$ws.Cells.Item(163, 15).Value = 18000
$ws.Cells.Item(163, 16).Value = 17500
$ws.Cells.Item(163, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(163, 18).Value = "Ecuador"
$ws.Cells.Item(163, 19).Value = 1250
$ws.Cells.Item(163, 20).Value = 14
